$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 ---
$ws.Range("J10").Value = 3
$ws.Range("M10").Value = ""

# --- Row 11 ---
$ws.Range("J11").Value = 0

# --- Row 12 ---
$ws.Range("M12").Value = ""

# --- Row 14 ---
$ws.Range("J14").Value = "Roblez"
$ws.Range("M14").Value = ""

# --- Row 16 ---
$ws.Range("J16").Value = "88-90 MPH"

# --- Row 17 ---
$ws.Range("J17").Value = "FB,CB,CH"

# --- Row 19 ---
$ws.Range("J19").Value = 4
$ws.Range("M19").Value = ""

# --- Row 21 ---
$ws.Range("M21").Value = ""

# --- Row 23 ---
$ws.Range("M23").Value = "Ground Ball"

# --- Row 24 ---
$ws.Range("M24").Value = "Single"

# --- Row 26 ---
$ws.Range("J26").Value = "FB,CB,CH"

# --- Row 28 ---
$ws.Range("J28").Value = 5
$ws.Range("M28").Value = ""

# --- Row 29 ---
$ws.Range("J29").Value = 1

# --- Row 30 ---
$ws.Range("M30").Value = ""

# --- Row 32 ---
$ws.Range("J32").Value = "Herbst"
$ws.Range("M32").Value = ""

# --- Row 34 ---
$ws.Range("J34").Value = "83-85 MPH"

# --- Row 35 ---
$ws.Range("J35").Value = "SL,FB,CB,CH"

# --- Row 37 ---
$ws.Range("J37").Value = 6
$ws.Range("M37").Value = ""

# --- Row 39 ---
$ws.Range("M39").Value = ""

# --- Row 41 ---
$ws.Range("J41").Value = "Herbst"
$ws.Range("M41").Value = "Line Drive"

# --- Row 42 ---
$ws.Range("J42").Value = "Right"
$ws.Range("M42").Value = "Out"

# --- Row 43 ---
$ws.Range("J43").Value = "83-85 MPH"

# --- Row 44 ---
$ws.Range("J44").Value = "SL,FB,CB,CH"

# --- Row 46 ---
$ws.Range("J46").Value = 7
$ws.Range("M46").Value = ""

# --- Row 47 ---
$ws.Range("J47").Value = 2

# --- Row 48 ---
$ws.Range("M48").Value = ""

# --- Row 50 ---
$ws.Range("J50").Value = "Plum"
$ws.Range("M50").Value = ""

# --- Row 51 ---
$ws.Range("M51").Value = "Undefined"

# --- Row 52 ---
$ws.Range("J52").Value = "84-86 MPH"

# --- Row 53 ---
$ws.Range("J53").Value = "SL,FB,CH"

# --- Row 61 ---
$ws.Range("J61").Value = 9
$ws.Range("M61").Value = ""

# --- Row 63 ---
$ws.Range("M63").Value = ""

# --- Row 65 ---
$ws.Range("J65").Value = "Thompson"
$ws.Range("M65").Value = "Ground Ball"

# --- Row 66 ---
$ws.Range("J66").Value = "Left"
$ws.Range("M66").Value = "Single"

# --- Row 67 ---
$ws.Range("J67").Value = "84-84 MPH"

# --- Row 68 ---
$ws.Range("J68").Value = "SL,FB,CH"
